$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.582.65"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.849.13"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").Value = "'1.029"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'321.01"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'1.024"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.4369"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "'0.07386"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.8814"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").Value = "'21.50"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.860.32"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "'5.506"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'6.703"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'0.07133"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "'85.01"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "'0.000009072"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'15.43"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "27.601.52"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'5.277"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "'11.22"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "2.090.50"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").Value = "'2.022"
$ws.Range("E25").Value = "  +5.69%  "
$ws.Range("D26").Value = "'157.08"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "'5.337"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Value = "'1.981"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "'117.48"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'0.08977"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").Value = "'0.7740"
$ws.Range("D33").Value = "'1.208"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'2.984"
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").Value = "'4.548"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").Value = "'1.024"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'1.138"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'0.01967"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'0.05249"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'2.858"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").Value = "'0.5167"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'0.1677"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "'6.822"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("D44").Value = "'8.801"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").Value = "'109.90"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "'10.66"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'0.06598"
$ws.Range("E47").Value = "  +4.19%  "
$ws.Range("D48").Value = "'1.025"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'1.698"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'0.4694"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "'1.888"
$ws.Range("E51").Value = "  -0.64%  "
